$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: cell reference, new text value.
# Values that parse as plain numbers get a leading apostrophe so Excel
# keeps them as text (matching the original inline-string cell contents)
# instead of silently converting them to numeric cells.
$updates = @(
    @("D2", "69.288.82"),
    @("E2", "  -0.59%  "),
    @("D3", "3.440.99"),
    @("E3", "  -1.79%  "),
    @("E4", "  -0.07%  "),
    @("D5", "'607.19"),
    @("E5", "  +0.39%  "),
    @("D6", "'167.21"),
    @("E6", "  -3.83%  "),
    @("D7", "3.434.91"),
    @("E7", "  -1.92%  "),
    @("D8", "'0.595"),
    @("E8", "  -2.05%  "),
    @("E9", "  +0.13%  "),
    @("E10", "  +0.79%  "),
    @("D11", "'7.02"),
    @("E11", "  -3.12%  "),
    @("D12", "'0.562"),
    @("E12", "  -3.26%  "),
    @("D13", "'44.16"),
    @("E13", "  -4.41%  "),
    @("E14", "  -1.71%  "),
    @("D15", "3.997.10"),
    @("E15", "  -1.73%  "),
    @("D16", "'8.14"),
    @("E16", "  -1.51%  "),
    @("B17", "WrappedEther"),
    @("C17", "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"),
    @("D17", "3.443.86"),
    @("E17", "  -1.71%  "),
    @("B18", "BitcoinCash"),
    @("C18", "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"),
    @("D18", "'580.50"),
    @("E18", "  -4.42%  "),
    @("B19", "WrappedBTC"),
    @("C19", "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"),
    @("D19", "69.330.88"),
    @("E19", "  -0.61%  "),
    @("D20", "'0.120"),
    @("E20", "  +0.93%  "),
    @("D21", "'17.09"),
    @("E21", "  -0.26%  "),
    @("D22", "'0.843"),
    @("E22", "  -3.06%  "),
    @("E23", "  -2.10%  "),
    @("D24", "'96.22"),
    @("E24", "  +0.67%  "),
    @("D25", "'15.12"),
    @("E25", "  -2.09%  "),
    @("D26", "'3.60"),
    @("E26", "  -2.92%  "),
    @("E27", "  +0.03%  "),
    @("E28", "  -5.24%  "),
    @("E29", "  -3.45%  "),
    @("D30", "'8.62"),
    @("E30", "  -3.81%  "),
    @("D31", "'7.80"),
    @("E31", "  -3.48%  "),
    @("D32", "'2.79"),
    @("E32", "  -6.33%  "),
    @("E33", "  -2.87%  "),
    @("D34", "'6.55"),
    @("E34", "  -4.97%  "),
    @("D35", "'579.95"),
    @("E35", "  -16.17%  "),
    @("D36", "'10.50"),
    @("E36", "  -1.66%  "),
    @("E37", "  -0.55%  "),
    @("D38", "'0.0953"),
    @("E38", "  -4.53%  "),
    @("E39", "  +0.29%  "),
    @("D40", "'56.02"),
    @("E40", "  -0.69%  "),
    @("E41", "  -0.64%  "),
    @("E42", "  -10.76%  "),
    @("D43", "3.233.94"),
    @("E43", "  -2.52%  "),
    @("D44", "0.0₃0684"),
    @("E44", "  -0.52%  "),
    @("B45", "InjectiveProtocol"),
    @("C45", "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"),
    @("D45", "'31.05"),
    @("E45", "  -3.61%  "),
    @("B46", "TheGraph"),
    @("C46", "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"),
    @("D46", "'0.295"),
    @("E46", "  -5.45%  "),
    @("E47", "  -5.32%  "),
    @("D48", "'2.39"),
    @("E48", "  -6.02%  "),
    @("E49", "  -2.93%  "),
    @("D50", "'134.28"),
    @("E50", "  +0.70%  "),
)

foreach ($u in $updates) {
    $ws.Range($u[0]).Value = $u[1]
}

Write-Output "Applied $($updates.Count) cell updates"
